# Add three new rows (3, 4, 5) of mile-post data to the "parsed mile posts" sheet,
# mirroring the existing row 2 pattern, then move the selection to C10 (matching
# the author's final cursor position recorded in the saved workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parsed mile posts")

# Row 3
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 10.15
$ws.Range("C3").Value = 11.56
$ws.Range("D3").Value = 120000
$ws.Range("E3").Value = "IS"
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

# Row 4
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 13.42
$ws.Range("C4").Value = 13.78
$ws.Range("D4").Value = 141000
$ws.Range("E4").Value = "IS"
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 3

# Row 5
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 10.15
$ws.Range("C5").Value = 11.56
$ws.Range("D5").Value = 120000
$ws.Range("E5").Value = "IS"
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 2

# Match the author's final selection position recorded in the file
$ws.Range("C10").Select()
